$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a = $ws.Cells.Item(150, 1)
$a.Formula = "=""05-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(150, 2).Value = 3412
$ws.Cells.Item(150, 3).Value = 623
$ws.Cells.Item(150, 4).Value = 505
$ws.Cells.Item(150, 5).Value = 647
$ws.Cells.Item(150, 6).Value = 712
$ws.Cells.Item(150, 7).Value = 925

$a = $ws.Cells.Item(151, 1)
$a.Formula = "=""06-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(151, 2).Value = 3745
$ws.Cells.Item(151, 3).Value = 902
$ws.Cells.Item(151, 4).Value = 615
$ws.Cells.Item(151, 5).Value = 880
$ws.Cells.Item(151, 6).Value = 546
$ws.Cells.Item(151, 7).Value = 802

$a = $ws.Cells.Item(152, 1)
$a.Formula = "=""09-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(152, 2).Value = 4766
$ws.Cells.Item(152, 3).Value = 1078
$ws.Cells.Item(152, 4).Value = 813
$ws.Cells.Item(152, 5).Value = 517
$ws.Cells.Item(152, 6).Value = 1158
$ws.Cells.Item(152, 7).Value = 1199

$a = $ws.Cells.Item(153, 1)
$a.Formula = "=""10-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(153, 2).Value = 4395
$ws.Cells.Item(153, 3).Value = 541
$ws.Cells.Item(153, 4).Value = 637
$ws.Cells.Item(153, 5).Value = 1283
$ws.Cells.Item(153, 6).Value = 983
$ws.Cells.Item(153, 7).Value = 951

$a = $ws.Cells.Item(154, 1)
$a.Formula = "=""11-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(154, 2).Value = 3835
$ws.Cells.Item(154, 3).Value = 455
$ws.Cells.Item(154, 4).Value = 476
$ws.Cells.Item(154, 5).Value = 1012
$ws.Cells.Item(154, 6).Value = 1031
$ws.Cells.Item(154, 7).Value = 861

$a = $ws.Cells.Item(155, 1)
$a.Formula = "=""12-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(155, 2).Value = 6171
$ws.Cells.Item(155, 3).Value = 884
$ws.Cells.Item(155, 4).Value = 952
$ws.Cells.Item(155, 5).Value = 1556
$ws.Cells.Item(155, 6).Value = 1356
$ws.Cells.Item(155, 7).Value = 1422

$a = $ws.Cells.Item(156, 1)
$a.Formula = "=""13-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(156, 2).Value = 4590
$ws.Cells.Item(156, 3).Value = 502
$ws.Cells.Item(156, 4).Value = 761
$ws.Cells.Item(156, 5).Value = 869
$ws.Cells.Item(156, 6).Value = 1344
$ws.Cells.Item(156, 7).Value = 1113

$a = $ws.Cells.Item(157, 1)
$a.Formula = "=""16-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(157, 2).Value = 3723
$ws.Cells.Item(157, 3).Value = 581
$ws.Cells.Item(157, 4).Value = 508
$ws.Cells.Item(157, 5).Value = 892
$ws.Cells.Item(157, 6).Value = 903
$ws.Cells.Item(157, 7).Value = 839

$a = $ws.Cells.Item(158, 1)
$a.Formula = "=""17-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(158, 2).Value = 5438
$ws.Cells.Item(158, 3).Value = 634
$ws.Cells.Item(158, 4).Value = 775
$ws.Cells.Item(158, 5).Value = 1335
$ws.Cells.Item(158, 6).Value = 1341
$ws.Cells.Item(158, 7).Value = 1352

$a = $ws.Cells.Item(159, 1)
$a.Formula = "=""18-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(159, 2).Value = 3907
$ws.Cells.Item(159, 3).Value = 530
$ws.Cells.Item(159, 4).Value = 614
$ws.Cells.Item(159, 5).Value = 685
$ws.Cells.Item(159, 6).Value = 1107
$ws.Cells.Item(159, 7).Value = 970

$a = $ws.Cells.Item(160, 1)
$a.Formula = "=""19-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(160, 2).Value = 4984
$ws.Cells.Item(160, 3).Value = 1505
$ws.Cells.Item(160, 4).Value = 645
$ws.Cells.Item(160, 5).Value = 1259
$ws.Cells.Item(160, 6).Value = 792
$ws.Cells.Item(160, 7).Value = 784

$a = $ws.Cells.Item(161, 1)
$a.Formula = "=""20-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(161, 2).Value = 4399
$ws.Cells.Item(161, 3).Value = 300
$ws.Cells.Item(161, 4).Value = 518
$ws.Cells.Item(161, 5).Value = 1540
$ws.Cells.Item(161, 6).Value = 1028
$ws.Cells.Item(161, 7).Value = 1013

$a = $ws.Cells.Item(162, 1)
$a.Formula = "=""23-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(162, 2).Value = 3208
$ws.Cells.Item(162, 3).Value = 614
$ws.Cells.Item(162, 4).Value = 953
$ws.Cells.Item(162, 5).Value = 510
$ws.Cells.Item(162, 6).Value = 513
$ws.Cells.Item(162, 7).Value = 617

$a = $ws.Cells.Item(163, 1)
$a.Formula = "=""24-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(163, 2).Value = 5374
$ws.Cells.Item(163, 3).Value = 542
$ws.Cells.Item(163, 4).Value = 622
$ws.Cells.Item(163, 5).Value = 1441
$ws.Cells.Item(163, 6).Value = 1448
$ws.Cells.Item(163, 7).Value = 1321

$a = $ws.Cells.Item(164, 1)
$a.Formula = "=""25-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(164, 2).Value = 4170
$ws.Cells.Item(164, 3).Value = 657
$ws.Cells.Item(164, 4).Value = 697
$ws.Cells.Item(164, 5).Value = 615
$ws.Cells.Item(164, 6).Value = 1072
$ws.Cells.Item(164, 7).Value = 1129

$a = $ws.Cells.Item(165, 1)
$a.Formula = "=""26-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(165, 2).Value = 3924
$ws.Cells.Item(165, 3).Value = 646
$ws.Cells.Item(165, 4).Value = 760
$ws.Cells.Item(165, 5).Value = 1012
$ws.Cells.Item(165, 6).Value = 776
$ws.Cells.Item(165, 7).Value = 729

$a = $ws.Cells.Item(166, 1)
$a.Formula = "=""27-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(166, 2).Value = 3930
$ws.Cells.Item(166, 3).Value = 581
$ws.Cells.Item(166, 4).Value = 836
$ws.Cells.Item(166, 5).Value = 754
$ws.Cells.Item(166, 6).Value = 976
$ws.Cells.Item(166, 7).Value = 783

$a = $ws.Cells.Item(167, 1)
$a.Formula = "=""30-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(167, 2).Value = 3051
$ws.Cells.Item(167, 3).Value = 321
$ws.Cells.Item(167, 4).Value = 481
$ws.Cells.Item(167, 5).Value = 407
$ws.Cells.Item(167, 6).Value = 1039
$ws.Cells.Item(167, 7).Value = 804

$a = $ws.Cells.Item(168, 1)
$a.Formula = "=""31-08-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(168, 2).Value = 4117
$ws.Cells.Item(168, 3).Value = 411
$ws.Cells.Item(168, 4).Value = 825
$ws.Cells.Item(168, 5).Value = 598
$ws.Cells.Item(168, 6).Value = 1359
$ws.Cells.Item(168, 7).Value = 923

$a = $ws.Cells.Item(169, 1)
$a.Formula = "=""01-09-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(169, 2).Value = 6431
$ws.Cells.Item(169, 3).Value = 524
$ws.Cells.Item(169, 4).Value = 1147
$ws.Cells.Item(169, 5).Value = 1404
$ws.Cells.Item(169, 6).Value = 2222
$ws.Cells.Item(169, 7).Value = 1134

$a = $ws.Cells.Item(170, 1)
$a.Formula = "=""02-09-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(170, 2).Value = 6203
$ws.Cells.Item(170, 3).Value = 884
$ws.Cells.Item(170, 4).Value = 989
$ws.Cells.Item(170, 5).Value = 995
$ws.Cells.Item(170, 6).Value = 1879
$ws.Cells.Item(170, 7).Value = 1456

$a = $ws.Cells.Item(171, 1)
$a.Formula = "=""03-09-2021"""
$a.Copy()
$a.PasteSpecial(-4163)
$ws.Cells.Item(171, 2).Value = 5020
$ws.Cells.Item(171, 3).Value = 187
$ws.Cells.Item(171, 4).Value = 492
$ws.Cells.Item(171, 5).Value = 920
$ws.Cells.Item(171, 6).Value = 1959
$ws.Cells.Item(171, 7).Value = 1463

$excel.CutCopyMode = 0
Write-Host "done"
